# Update the "main pages" (Metadata sheet) of the ValueSet workbook:
#  - Translate the Publisher and Contact values from German to English
#  - Add the missing Description text

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"
$ws.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"
$ws.Range("B12").Value = "consent states - minimal subset CONSENT documents"
